$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# A new FMOD event (":/Drums/StickCollison") was inserted as row 4, pushing
# the existing rows 4-14 down to rows 5-15. Emulate the row-insert by
# shifting the rows downward one at a time, bottom-up, copying both the
# formatting and the values (this engine's PasteSpecial needs formats and
# values pasted separately, and the destination cleared first, otherwise
# stale values/styles are left behind when the source cell is blank).
# ---------------------------------------------------------------------------
for ($r = 14; $r -ge 4; $r--) {
    $src = $ws.Range("A" + $r + ":G" + $r)
    $dst = $ws.Range("A" + ($r + 1) + ":G" + ($r + 1))
    $dst.ClearContents()
    $src.Copy()
    $dst.PasteSpecial(-4122) | Out-Null   # xlPasteFormats
    $dst.PasteSpecial(-4163) | Out-Null   # xlPasteValues
}
$excel.CutCopyMode = 0

# Fill in the newly-freed row 4 with the new event's data.
$ws.Range("A4").Value = ":/Drums/StickCollison"
$ws.Range("B4").Value = "N"
$ws.Range("C4").Value = "N"
$ws.Range("D4").Value = "StickCollision"
$ws.Range("E4").Value = "Discrete"
$ws.Range("F4").Value = "0-1"

# ---------------------------------------------------------------------------
# A second new event (":/Atmos/Haunting") was appended as the new last row
# (row 16). Copy the formatting from the previous last row (15) first, then
# set its values.
# ---------------------------------------------------------------------------
$lastSrc = $ws.Range("A15:G15")
$newRow = $ws.Range("A16:G16")
$lastSrc.Copy()
$newRow.PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("A16").Value = ":/Atmos/Haunting"
$ws.Range("B16").Value = "Y"
$ws.Range("C16").Value = "N"

$ws.Range("F4").Select()
